$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12482.267
$ws.Range("I6").Value = 13956.385
$ws.Range("K6").Value = 41869.155
$ws.Range("M6").Value = -41757.155
$ws.Range("H33").Value = 1351635.2
$ws.Range("I33").Value = 2252533.5
$ws.Range("J33").Value = 288
$ws.Range("K33").Value = 2252533.5
$ws.Range("L33").Value = 288
$ws.Range("M33").Value = -2252304.5
$ws.Range("N33").Value = -746
$ws.Range("H38").Value = 695.75
$ws.Range("I38").Value = 427.66666
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 1282.99998
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = -910.9999800000001
$ws.Range("N38").Value = -5244
$ws.Range("H39").Value = 195
$ws.Range("I39").Value = 195
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 585
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -289
$ws.Range("H107").Value = 5834.6665
$ws.Range("I107").Value = 4750
$ws.Range("K107").Value = 4750
$ws.Range("M107").Value = -2830
$ws.Range("H111").Value = 1057
$ws.Range("I111").Value = 1057
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3171
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -104
$ws.Range("H115").Value = 1560
$ws.Range("I115").Value = 933.3333
$ws.Range("J115").Value = 2500
$ws.Range("K115").Value = 2799.9999
$ws.Range("L115").Value = 7500
$ws.Range("M115").Value = -1232.9999
$ws.Range("N115").Value = -10634
$ws.Range("H131").Value = 6220.9
$ws.Range("I131").Value = 6165.875
$ws.Range("J131").Value = 6441
$ws.Range("K131").Value = 18497.625
$ws.Range("L131").Value = 19323
$ws.Range("M131").Value = -13457.625
$ws.Range("N131").Value = -29403
$ws.Range("H138").Value = 4701.0234
$ws.Range("I138").Value = 3451.8667
$ws.Range("J138").Value = 5370.2144
$ws.Range("K138").Value = 10355.6001
$ws.Range("L138").Value = 16110.6432
$ws.Range("M138").Value = -5215.6001
$ws.Range("N138").Value = -26390.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2993.796
$ws.Range("I32").Value = 1694.3334
$ws.Range("K32").Value = 1694.3334
$ws.Range("M32").Value = -1407.3334
$ws.Range("H61").Value = 7656.5356
$ws.Range("I61").Value = 5636.625
$ws.Range("K61").Value = 5636.625
$ws.Range("M61").Value = -5424.625
$ws.Range("H101").Value = 39000
$ws.Range("J101").Value = 39000
$ws.Range("L101").Value = 39000
$ws.Range("N101").Value = -45490
$ws.Range("H110").Value = 1518.1923
$ws.Range("I110").Value = 1452.3182
$ws.Range("J110").Value = 1880.5
$ws.Range("K110").Value = 1452.3182
$ws.Range("L110").Value = 1880.5
$ws.Range("M110").Value = 592.6818000000001
$ws.Range("N110").Value = -5970.5
$ws.Range("H124").Value = 53065.4
$ws.Range("J124").Value = 53065.4
$ws.Range("L124").Value = 53065.4
$ws.Range("N124").Value = -62885.4
$ws.Range("H125").Value = 69000
$ws.Range("J125").Value = 69000
$ws.Range("L125").Value = 69000
$ws.Range("N125").Value = -78840
$ws.Range("H132").Value = 7356.2183
$ws.Range("I132").Value = 6120.5435
$ws.Range("K132").Value = 18361.6305
$ws.Range("M132").Value = -15831.6305
$ws.Range("H136").Value = 7656.5356
$ws.Range("I136").Value = 5636.625
$ws.Range("K136").Value = 16909.875
$ws.Range("M136").Value = -14359.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1931.8334
$ws.Range("I99").Value = 1931.8334
$ws.Range("K99").Value = 1931.8334
$ws.Range("M99").Value = -433.8334
$ws.Range("H107").Value = 1878.8334
$ws.Range("I107").Value = 1963.8334
$ws.Range("J107").Value = 1793.8334
$ws.Range("K107").Value = 1963.8334
$ws.Range("L107").Value = 1793.8334
$ws.Range("M107").Value = -43.83339999999998
$ws.Range("N107").Value = -5633.8334
$ws.Range("H134").Value = 10878.429
$ws.Range("I134").Value = 8288.333000000001
$ws.Range("K134").Value = 24864.999
$ws.Range("M134").Value = -22329.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 520.5625
$ws.Range("I22").Value = 335.7
$ws.Range("J22").Value = 828.6667
$ws.Range("K22").Value = 335.7
$ws.Range("L22").Value = 828.6667
$ws.Range("M22").Value = 14.30000000000001
$ws.Range("N22").Value = -1528.6667
$ws.Range("H31").Value = 6047.573
$ws.Range("I31").Value = 1475.4828
$ws.Range("K31").Value = 1475.4828
$ws.Range("M31").Value = -1180.4828
$ws.Range("H34").Value = 6047.573
$ws.Range("I34").Value = 1475.4828
$ws.Range("K34").Value = 1475.4828
$ws.Range("M34").Value = -1273.4828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 356.15
$ws.Range("I2").Value = 142.4
$ws.Range("J2").Value = 569.9
$ws.Range("K2").Value = 854.4000000000001
$ws.Range("L2").Value = 3419.4
$ws.Range("M2").Value = -741.4000000000001
$ws.Range("N2").Value = -3645.4
$ws.Range("H7").Value = 100000080
$ws.Range("I7").Value = 100000080
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 300000240
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -300000128
$ws.Range("H11").Value = 209.76923
$ws.Range("I11").Value = 125.8
$ws.Range("J11").Value = 489.66666
$ws.Range("K11").Value = 377.4
$ws.Range("L11").Value = 1468.99998
$ws.Range("M11").Value = -237.4
$ws.Range("N11").Value = -1748.99998
$ws.Range("H34").Value = 4243
$ws.Range("I34").Value = 246.75
$ws.Range("J34").Value = 5575.0835
$ws.Range("K34").Value = 740.25
$ws.Range("L34").Value = 16725.2505
$ws.Range("M34").Value = -656.25
$ws.Range("N34").Value = -16893.2505
$ws.Range("H39").Value = 4624.5
$ws.Range("J39").Value = 5332.6665
$ws.Range("L39").Value = 15997.9995
$ws.Range("N39").Value = -16585.9995
$ws.Range("H46").Value = 11500408
$ws.Range("J46").Value = 318
$ws.Range("L46").Value = 954
$ws.Range("N46").Value = -1136
$ws.Range("H50").Value = 1338.174
$ws.Range("J50").Value = 1427.5714
$ws.Range("L50").Value = 4282.7142
$ws.Range("N50").Value = -5244.7142
$ws.Range("H53").Value = 1338.174
$ws.Range("J53").Value = 1427.5714
$ws.Range("L53").Value = 4282.7142
$ws.Range("N53").Value = -5244.7142
$ws.Range("H68").Value = 2767.0417
$ws.Range("I68").Value = 2736.625
$ws.Range("J68").Value = 2782.25
$ws.Range("K68").Value = 8209.875
$ws.Range("L68").Value = 8346.75
$ws.Range("M68").Value = -7398.875
$ws.Range("N68").Value = -9968.75
$ws.Range("H71").Value = 2767.0417
$ws.Range("I71").Value = 2736.625
$ws.Range("J71").Value = 2782.25
$ws.Range("K71").Value = 24629.625
$ws.Range("L71").Value = 25040.25
$ws.Range("M71").Value = -20573.625
$ws.Range("N71").Value = -33152.25
$ws.Range("H81").Value = 2555.5557
$ws.Range("J81").Value = 2555.5557
$ws.Range("L81").Value = 7666.6671
$ws.Range("N81").Value = -9912.667099999999
$ws.Range("H84").Value = 2555.5557
$ws.Range("J84").Value = 2555.5557
$ws.Range("L84").Value = 23000.0013
$ws.Range("N84").Value = -34232.0013
$ws.Range("H107").Value = 2543.2727
$ws.Range("I107").Value = 2508.6667
$ws.Range("J107").Value = 2699
$ws.Range("K107").Value = 7526.000100000001
$ws.Range("L107").Value = 8097
$ws.Range("M107").Value = -5606.000100000001
$ws.Range("N107").Value = -11937

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 15000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15340
$ws.Range("H43").Value = 468384.62
$ws.Range("I43").Value = 450000
$ws.Range("J43").Value = 484142.84
$ws.Range("K43").Value = 450000
$ws.Range("L43").Value = 484142.84
$ws.Range("M43").Value = -449807
$ws.Range("N43").Value = -484528.84
$ws.Range("H93").Value = 1911.2
$ws.Range("I93").Value = 1940.9744
$ws.Range("J93").Value = 750
$ws.Range("K93").Value = 1940.9744
$ws.Range("L93").Value = 750
$ws.Range("M93").Value = -692.9744000000001
$ws.Range("N93").Value = -3246
$ws.Range("H100").Value = 509625.72
$ws.Range("I100").Value = 4981.8945
$ws.Range("J100").Value = 3705703.2
$ws.Range("K100").Value = 4981.8945
$ws.Range("L100").Value = 3705703.2
$ws.Range("M100").Value = -4440.8945
$ws.Range("N100").Value = -3706785.2
$ws.Range("H127").Value = 52100
$ws.Range("J127").Value = 52100
$ws.Range("L127").Value = 52100
$ws.Range("N127").Value = -62020
$ws.Range("H132").Value = 12350857
$ws.Range("I132").Value = 12350857
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 37052571
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -37050041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 17500
$ws.Range("I15").Value = 17500
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 17500
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -17212
$ws.Range("N15").ClearContents()
$ws.Range("H69").Value = 58000
$ws.Range("J69").Value = 58000
$ws.Range("L69").Value = 58000
$ws.Range("N69").Value = -59498
$ws.Range("H72").Value = 58000
$ws.Range("J72").Value = 58000
$ws.Range("L72").Value = 174000
$ws.Range("N72").Value = -181488
$ws.Range("H107").Value = 1136.8572
$ws.Range("I107").Value = 971.5294
$ws.Range("K107").Value = 2914.5882
$ws.Range("M107").Value = -994.5882000000001
$ws.Range("H132").Value = 19046.568
$ws.Range("I132").Value = 13443.511
$ws.Range("K132").Value = 40330.533
$ws.Range("M132").Value = -37800.533
$ws.Range("H136").Value = 2898.182
$ws.Range("I136").Value = 2897.889
$ws.Range("J136").Value = 2899.5
$ws.Range("K136").Value = 8693.667000000001
$ws.Range("L136").Value = 2899.5
$ws.Range("M136").Value = -6143.667000000001
$ws.Range("N136").Value = -13798.5
